$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: 2016-08-29, matching the date styling used by the rows above it
# (copy A7's format so the new date cell reuses the existing date style
# instead of Excel minting a brand-new numFmt/cellXf).
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A8").Value = 42611

# All tracked items are complete for this row.
$ws.Range("B8:L8").Value = "Done"

$ws.Range("L8").Select()
